# Generate Report for Handback
# Updates the localization-status workbook: marks zh-cn/de-de as handed back,
# records the handback target/xlf file + timestamp, and links to the source doc.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$mdFileName = "ef652101-77ef-4c7c-ab92-81b405b701f1.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9eee9d4a28ad7a476fa1210c5e5546fd65260007/e2e/ef652101-77ef-4c7c-ab92-81b405b701f1.md"
$zhXlf      = "ef652101-77ef-4c7c-ab92-81b405b701f1.d0c044901a345deae0897cab2e9f8383a5c0b8b4.zh-cn.xlf"
$deXlf      = "ef652101-77ef-4c7c-ab92-81b405b701f1.d0c044901a345deae0897cab2e9f8383a5c0b8b4.de-de.xlf"

# --- Status column updates (Overview!E2/F2, zh-cn!C2, de-de!C2 all share text) ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText

# --- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Range("J2").Value = $zhXlf
$wsZhCn.Range("K2").Value = "2016-08-15 14:56:16"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276

# --- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Range("J2").Value = $deXlf
$wsDeDe.Range("K2").Value = "2016-08-15 14:56:23"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276

# --- Column width adjustments (status column + new target/handback columns widen) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.0   # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.0   # F

$wsZhCn.Columns.Item(3).ColumnWidth = 29.0    # C (Status)
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15   # I (Latest Target File)
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15  # J (Latest Handback File)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.0    # C (Status)
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15   # I (Latest Target File)
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15  # J (Latest Handback File)

Write-Host "Handback report generated"
